# PowerShell-style Excel COM-interop script
# Implements commit: "Atualiza dados da liga classica"
# Adds a new team "TEAM LOPES 99" with zeroed scores across all 9 sheets of the workbook,
# and re-populates the team-name column (A) of the six single-column standings sheets
# ("Classif Turno 2" and the five "Mes - ..." sheets) to match the edited source state.

$wb = $excel.ActiveWorkbook

$newTeam = "TEAM LOPES 99"

# ---------------------------------------------------------------------------
# "Geral" sheet (full round-by-round table, columns B:AM = Rodada 1..38)
# ---------------------------------------------------------------------------
$wsGeral = $wb.Worksheets.Item("Geral")
$wsGeral.Range("A52").Copy($wsGeral.Range("A53"))
$wsGeral.Range("A53").Value = $newTeam
for ($col = 2; $col -le 39; $col++) {
    $wsGeral.Cells.Item(53, $col).Value = 0
}

# ---------------------------------------------------------------------------
# "Turno 2" sheet (round-by-round table, columns B:T = Rodada 1..19)
# ---------------------------------------------------------------------------
$wsTurno2 = $wb.Worksheets.Item("Turno 2")
$wsTurno2.Range("A52").Copy($wsTurno2.Range("A53"))
$wsTurno2.Range("A53").Value = $newTeam
for ($col = 2; $col -le 20; $col++) {
    $wsTurno2.Cells.Item(53, $col).Value = 0
}

# ---------------------------------------------------------------------------
# Standings sheets: column A (team name) + column B (points, always 0).
# These six sheets share the exact same team ordering after the edit.
# ---------------------------------------------------------------------------
$standingsSheetNames = @(
    "Classif Turno 2",
    "Mês - Janeiro",
    "Mês - Fevereiro",
    "Mês - Março",
    "Mês - Abril",
    "Mês - Maio",
    "Mês - Julho"
)

$teamOrder = @(
    "A Lenda Super Vasco F.c",
    "A Lenda Super Vascão f.c",
    "mercearia Estrela",
    "Máquina Laranjja",
    "NaoVaiDescer!",
    "Paulo Virgili FC",
    "Pity10",
    "pra sempre imortal fc",
    "Profit Soccer",
    "PUXE FC",
    "Rolo Compressor ZN",
    "S.E.R. GRILLO",
    "seralex",
    "SERGRILLO",
    "Sport Clube PAIM",
    "Super Vasco f.c",
    "SUPER VASCÃO F.C",
    "Tabajara de Inhaua PB1",
    "TATITTA FC",
    "Tatols Beants F.C",
    "TEAM LOPES 99",
    "teves_futsal20 f.c",
    "Texas Club 2026",
    "TIGRE LEON",
    "Time do S.A.P.O",
    "Mau Humor F.C.",
    "MAFRA MARTINS FC",
    "Luis lemes inter",
    "FBC Colorado II",
    "AZURRA82",
    "Bandoleros FCS",
    "BordonFC04",
    "C.A. Charru@",
    "cartola scheuer17",
    "CARTOLEIRO DO VALLE PRO26.5",
    "dasdoresfc",
    "DM Studio",
    "Dom Camillo68",
    "FBC Colorado",
    "FC castelo Branco 2",
    "lsauer fc",
    "FC Los Castilho",
    "Fedato Futebol Clube",
    "FIGUEIRA DA ILHA",
    "FÚRIA LEON",
    "Gig@ntte",
    "Gremiomaniasm",
    "Grêmio imortal 37",
    "JUV. KP",
    "JV5 Tricolor Gaúcho",
    "LISI GREMISTA",
    "VASCO MARTINS FC"
)

foreach ($sheetName in $standingsSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Extend the used range (row 53) with the same formatting as row 52 first.
    $ws.Range("A52:B52").Copy($ws.Range("A53:B53"))

    for ($i = 0; $i -lt $teamOrder.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $teamOrder[$i]
        $ws.Cells.Item($row, 2).Value = 0
    }
}

